$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 106
$ws.Range("H106").Value = 9344602
$ws.Range("I106").Value = 12458670
$ws.Range("J106").Value = 2398.6667
$ws.Range("K106").Value = 12458670
$ws.Range("L106").Value = 2398.6667
$ws.Range("M106").Value = -12458039
$ws.Range("N106").Value = -3660.6667
# Row 112
$ws.Range("H112").Value = 10910181
$ws.Range("J112").Value = 10910181
$ws.Range("L112").Value = 32730543
$ws.Range("N112").Value = -32732759
# Row 113
$ws.Range("H113").Value = 175950.5
$ws.Range("I113").Value = 414282
$ws.Range("J113").Value = 5713.7144
$ws.Range("K113").Value = 414282
$ws.Range("L113").Value = 5713.7144
$ws.Range("M113").Value = -411028
$ws.Range("N113").Value = -12221.7144
# Row 116
$ws.Range("H116").Value = 8026399.5
$ws.Range("I116").Value = 12836438
$ws.Range("J116").Value = 9668.666999999999
$ws.Range("K116").Value = 12836438
$ws.Range("L116").Value = 9668.666999999999
$ws.Range("M116").Value = -12832996
$ws.Range("N116").Value = -16552.667
# Row 129
$ws.Range("H129").Value = 990.16
$ws.Range("I129").Value = 597
$ws.Range("J129").Value = 995.4729599999999
$ws.Range("K129").Value = 1791
$ws.Range("L129").Value = 2986.41888
$ws.Range("M129").Value = 3209
$ws.Range("N129").Value = -12986.41888
# Row 133
$ws.Range("H133").Value = 12666.667
$ws.Range("J133").Value = 12666.667
$ws.Range("L133").Value = 12666.667
$ws.Range("N133").Value = -22786.667

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4502.6523
$ws.Range("I32").Value = 2355.5334
$ws.Range("J32").Value = 8528.5
$ws.Range("K32").Value = 2355.5334
$ws.Range("L32").Value = 8528.5
$ws.Range("M32").Value = -2068.5334
$ws.Range("N32").Value = -9102.5
# Row 61
$ws.Range("H61").Value = 2550.0667
$ws.Range("I61").Value = 1828.174
$ws.Range("K61").Value = 1828.174
$ws.Range("M61").Value = -1616.174
# Row 74
$ws.Range("H74").Value = 10188.267
$ws.Range("I74").Value = 2365.818
$ws.Range("K74").Value = 2365.818
$ws.Range("M74").Value = -1491.818
# Row 77
$ws.Range("H77").Value = 10188.267
$ws.Range("I77").Value = 2365.818
$ws.Range("K77").Value = 11829.09
$ws.Range("M77").Value = -7461.09
# Row 110
$ws.Range("H110").Value = 1079.1666
$ws.Range("J110").Value = 1533.3334
$ws.Range("L110").Value = 1533.3334
$ws.Range("N110").Value = -5623.3334
# Row 122
$ws.Range("H122").Value = 2525.7856
$ws.Range("I122").Value = 2604.3635
$ws.Range("J122").Value = 2237.6667
$ws.Range("K122").Value = 7813.0905
$ws.Range("L122").Value = 6713.000100000001
$ws.Range("M122").Value = -5363.0905
$ws.Range("N122").Value = -11613.0001
# Row 133
$ws.Range("H133").Value = 31000
$ws.Range("J133").Value = 31000
$ws.Range("L133").Value = 31000
$ws.Range("N133").Value = -36060
# Row 136
$ws.Range("H136").Value = 2550.0667
$ws.Range("I136").Value = 1828.174
$ws.Range("K136").Value = 5484.522
$ws.Range("M136").Value = -2934.522
# Row 139
$ws.Range("H139").Value = 60715
$ws.Range("J139").Value = 60715
$ws.Range("L139").Value = 60715
$ws.Range("N139").Value = -70995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2051.2222
$ws.Range("I20").Value = 2057.3
$ws.Range("J20").Value = 2043.625
$ws.Range("K20").Value = 2057.3
$ws.Range("L20").Value = 2043.625
$ws.Range("M20").Value = -1810.3
$ws.Range("N20").Value = -2537.625
# Row 94
$ws.Range("H94").Value = 1795.6666
$ws.Range("I94").Value = 2023
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 2023
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -1572
$ws.Range("N94").Value = -1902
# Row 107
$ws.Range("H107").Value = 451.33334
$ws.Range("I107").Value = 310
$ws.Range("J107").Value = 538.3077
$ws.Range("K107").Value = 310
$ws.Range("L107").Value = 538.3077
$ws.Range("M107").Value = 1610
$ws.Range("N107").Value = -4378.3077

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 782
$ws.Range("I16").Value = 470.33334
$ws.Range("J16").Value = 1249.5
$ws.Range("K16").Value = 470.33334
$ws.Range("L16").Value = 1249.5
$ws.Range("M16").Value = -183.33334
$ws.Range("N16").Value = -1823.5
# Row 58
$ws.Range("H58").Value = 2189.6086
$ws.Range("I58").Value = 1464.8
$ws.Range("K58").Value = 1464.8
$ws.Range("M58").Value = -1261.8
# Row 113
$ws.Range("H113").Value = 782
$ws.Range("I113").Value = 470.33334
$ws.Range("J113").Value = 1249.5
$ws.Range("K113").Value = 470.33334
$ws.Range("L113").Value = 1249.5
$ws.Range("M113").Value = 1699.66666
$ws.Range("N113").Value = -5589.5
# Row 136
$ws.Range("H136").Value = 2189.6086
$ws.Range("I136").Value = 1464.8
$ws.Range("K136").Value = 4394.4
$ws.Range("M136").Value = -1844.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 103
$ws.Range("H103").Value = 1590.8055
$ws.Range("I103").Value = 504.46155
$ws.Range("J103").Value = 2204.8262
$ws.Range("K103").Value = 1513.38465
$ws.Range("L103").Value = 6614.4786
$ws.Range("M103").Value = -634.38465
$ws.Range("N103").Value = -8372.4786
# Row 129
$ws.Range("H129").Value = 2110.6924
$ws.Range("I129").Value = 1372
$ws.Range("J129").Value = 2572.375
$ws.Range("K129").Value = 4116
$ws.Range("L129").Value = 7717.125
$ws.Range("M129").Value = 884
$ws.Range("N129").Value = -17717.125
# Row 131
$ws.Range("H131").Value = 1778.6578
$ws.Range("I131").Value = 623.3333
$ws.Range("J131").Value = 1877.6857
$ws.Range("K131").Value = 1869.9999
$ws.Range("L131").Value = 5633.0571
$ws.Range("M131").Value = 3170.0001
$ws.Range("N131").Value = -15713.0571
# Row 132
$ws.Range("H132").Value = 1422.1428
# Row 133
$ws.Range("H133").Value = 6643
$ws.Range("J133").Value = 11780
$ws.Range("L133").Value = 35340
$ws.Range("N133").Value = -45460

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 161.36
$ws.Range("I107").Value = 152.2
$ws.Range("J107").Value = 175.1
$ws.Range("K107").Value = 152.2
$ws.Range("L107").Value = 175.1
$ws.Range("M107").Value = 1767.8
$ws.Range("N107").Value = -4015.1
# Row 132
$ws.Range("H132").Value = 2439.5454
$ws.Range("I132").Value = 1504.5
$ws.Range("K132").Value = 4513.5
$ws.Range("M132").Value = -1983.5
# Row 138
$ws.Range("H138").Value = 60500
$ws.Range("J138").Value = 60500
$ws.Range("L138").Value = 60500
$ws.Range("N138").Value = -70780

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 623.4815
$ws.Range("I93").Value = 606.2727
$ws.Range("K93").Value = 606.2727
$ws.Range("M93").Value = 641.7273
# Row 132
$ws.Range("H132").Value = 3576.5518
$ws.Range("I132").Value = 2465.0588
$ws.Range("J132").Value = 5151.1665
$ws.Range("K132").Value = 7395.176399999999
$ws.Range("L132").Value = 15453.4995
$ws.Range("M132").Value = -4865.176399999999
$ws.Range("N132").Value = -20513.4995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 9166.666999999999
$ws.Range("J15").Value = 9166.666999999999
$ws.Range("L15").Value = 9166.666999999999
$ws.Range("N15").Value = -9742.666999999999
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()
# Row 62
$ws.Range("H62").Value = 14868.909
$ws.Range("I62").Value = 17708.75
$ws.Range("K62").Value = 17708.75
$ws.Range("M62").Value = -17084.75
# Row 65
$ws.Range("H65").Value = 14868.909
$ws.Range("I65").Value = 17708.75
$ws.Range("K65").Value = 88543.75
$ws.Range("M65").Value = -85423.75
# Row 122
$ws.Range("H122").Value = 60765.234
$ws.Range("I122").Value = 92573.45
$ws.Range("J122").Value = 2450.1667
$ws.Range("K122").Value = 277720.35
$ws.Range("L122").Value = 7350.500100000001
$ws.Range("M122").Value = -275270.35
$ws.Range("N122").Value = -12250.5001
# Row 132
$ws.Range("H132").Value = 11113639
$ws.Range("I132").Value = 13890812
$ws.Range("K132").Value = 41672436
$ws.Range("M132").Value = -41669906

Write-Host "Update complete"
